$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 100.00006
$ws.Range("C3").Value = 0.099999712
$ws.Range("D3").Value = 0.0009999997600000001
$ws.Range("E3").Value = 0.7102666230633887
$ws.Range("F3").Value = 1.034145083016212
$ws.Range("G3").Value = 80607.15869937603

# Row 4
$ws.Range("B4").Value = 17.33964079630801
$ws.Range("C4").Value = 0.01733953183537289
$ws.Range("D4").Value = 0.0001733934224350253
$ws.Range("E4").Value = 0.04147410335169903
$ws.Range("F4").Value = 0.003436597065559849
$ws.Range("G4").Value = 13407.49611982769

# Row 5
$ws.Range("B5").Value = 70
$ws.Range("C5").Value = 0.07000000000000001
$ws.Range("D5").Value = 0.0007
$ws.Range("E5").Value = 0.578895168070878
$ws.Range("F5").Value = 1.029532856100207
$ws.Range("G5").Value = 55989.3350985158

# Row 6
$ws.Range("B6").Value = 85
$ws.Range("C6").Value = 0.08500000000000001
$ws.Range("D6").Value = 0.00085
$ws.Range("E6").Value = 0.6826635816458004
$ws.Range("F6").Value = 1.031562654380683
$ws.Range("G6").Value = 69582.03205037788

# Row 7
$ws.Range("B7").Value = 100
$ws.Range("C7").Value = 0.1
$ws.Range("D7").Value = 0.001
$ws.Range("E7").Value = 0.7098506874547897
$ws.Range("F7").Value = 1.033217651792513
$ws.Range("G7").Value = 78792.86996330289

# Row 8
$ws.Range("B8").Value = 115
$ws.Range("C8").Value = 0.115
$ws.Range("D8").Value = 0.00115
$ws.Range("E8").Value = 0.7404884652144335
$ws.Range("F8").Value = 1.035783408489331
$ws.Range("G8").Value = 90208.41744320151

# Row 9
$ws.Range("B9").Value = 130
$ws.Range("C9").Value = 0.13
$ws.Range("D9").Value = 0.0013
$ws.Range("E9").Value = 0.8099847382030015
$ws.Range("F9").Value = 1.049276336254016
$ws.Range("G9").Value = 124214.6968227415
